$d = $word.ActiveDocument

function FindReplaceOne($findText, $replaceText) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)
    if (-not $ok) {
        Write-Host "WARNING: not found -> $findText"
    }
    return $ok
}

function FindOnly($findText) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Host "WARNING: not found -> $findText"
    }
    return $r
}

# ---------------------------------------------------------------
# 1. Personal Statement paragraph: replace entire paragraph text
# ---------------------------------------------------------------
$oldStatement = "For this study, I will serve as the principal investigator and work with the PEDSnet learning network, the study staff, and co-investigators to oversee all aspects of the research project, ensure the scientific integrity of the proposed study, and ensure its successful completion. My expertise and experience in leveraging high-resolution spatiotemporal exposure assessments in electronic health record to study psychiatric outcomes will contribute to the success of this proposed research. My personal hope is that this work will provide insight into biological mechanisms that mediate environmental risk to help identify modifiable primary interventions to reduce psychiatric morbidity in children and adolescents."
$newStatement = "As Carson's primary mentor for this training award, I am committed to providing a rigorous and supportive training environment that will equip him with the skills needed for an independent research career. His proposed training plan is strategically designed to integrate advanced computational methods, interdisciplinary collaboration, and translational research, ensuring he develops expertise in biomedical informatics and environmental epidemiology. Through structured mentorship, hands-on experience, and engagement with clinical and policy stakeholders, Carson will gain the technical proficiency, leadership skills, and professional network necessary to secure external funding and establish his own research program. I am confident that this training plan, combined with Carson's strong analytical abilities and scientific curiosity, will position him for success as a leader in data-driven public health research."
[void](FindReplaceOne $oldStatement $newStatement)

# ---------------------------------------------------------------
# 2. Add a new FirstParagraph-styled paragraph right after
#    "Peer-reviewed publications I would like to highlight include:"
# ---------------------------------------------------------------
$r = FindOnly "Peer-reviewed publications I would like to highlight include:"
$para = $r.Paragraphs(1)
[void]$para.Range.InsertParagraphAfter()
# re-find to get a fresh reference to the heading paragraph, then grab its successor
$r2 = FindOnly "Peer-reviewed publications I would like to highlight include:"
$headingPara = $r2.Paragraphs(1)
$newPara = $headingPara.Next()
$newPara.Range.Text = "These publications were selected to demonstrate my scientific productivity when I am the primary mentor of a clinician-scientist in training:"
$newPara.Style = "FirstParagraph"

# ---------------------------------------------------------------
# 3. Publication 1 (Milan N Parikh ... American Journal of Epidemiology)
# ---------------------------------------------------------------
[void](FindReplaceOne "Erika Rasnick Manning, Qing Duan, Stuart Taylor, Sarah Ray, Alexandra MS Corley, Joseph Michael, Ryan Gillette, Ndidi Unaka, David Hartley, Andrew F Beck, " "Milan N Parikh, Erika Rasnick Manning, Liang Niu, Anna Kotsakis Ruehlmann, Alonzo T Folger, Kelly J Brunst, ")
[void](FindReplaceOne ". Development of a Multimodal Geomarker Pipeline to Assess the Impact of Social, Economic, and Environmental Factors on Pediatric Health Outcomes. " ". Increasing Temporal Sensitivity of Omics Association Studies with Epigenome-Wide Distributed Lag Models. ")
[void](FindReplaceOne "Journal of the American Medical Informatics Association" "American Journal of Epidemiology")
[void](FindReplaceOne ". In press. 2024." ". In Press. 2024.    ")

# ---------------------------------------------------------------
# 4. Publication 2 (Harsimran Makkad ... Journal of Racial and Ethnic Health Disparities)
#    First prepend a new plain-text leading-author run before "Cole Brokamp"
# ---------------------------------------------------------------
$r3 = FindOnly ", Jeffrey R. Strawn, Andrew F. Beck, Pat Ryan. Pediatric Psychiatric Emergency Department Utilization and Fine Particulate Matter: A Case-Crossover Study. "
$para3 = $r3.Paragraphs(1)
$paraStart3 = $para3.Range.Start
$insertRange3 = $d.Range($paraStart3, $paraStart3)
$insertRange3.InsertBefore("Harsimran Makkad, Amisha Saini, Erika Rasnick Manning, Qing Duan, Stephen Colegate, ")

[void](FindReplaceOne ", Jeffrey R. Strawn, Andrew F. Beck, Pat Ryan. Pediatric Psychiatric Emergency Department Utilization and Fine Particulate Matter: A Case-Crossover Study. " ". Racial Fairness of Individual- and Community-Level Proxies of Socioeconomic Status Among Birthing Parent–Child Dyads. ")
[void](FindReplaceOne "Environmental Health Perspectives" "Journal of Racial and Ethnic Health Disparities")
[void](FindReplaceOne ". 127(9). 2019." ". Online. 2024.    ")

# ---------------------------------------------------------------
# 5. Publication 3 (Stephen P Colegate ... Journal of Clinical and Translational Science)
# ---------------------------------------------------------------
[void](FindReplaceOne "Erika Manning, Qing Duan, " "Stephen P Colegate, Anushka Palipana, Emrah Gecili, Rhonda D Szczesniak, ")
[void](FindReplaceOne ". Incorporating Parcel-Based Housing Conditions to Increase the Precision of Identifying Children with Elevated Blood Lead. " ". Evaluating Precision Medicine Tools in Cystic Fibrosis for Racial and Ethnic Fairness. ")
[void](FindReplaceOne "Journal of Public Health Management & Practice" "Journal of Clinical and Translational Science")
[void](FindReplaceOne ". In Press. 2024." ". In press. 2024.    ")

# ---------------------------------------------------------------
# 6. Publication 4 (Jordan Pennington ... American Journal of Health Promotion)
#    Prepend a new plain-text leading-author run before "Cole Brokamp"
# ---------------------------------------------------------------
$r4 = FindOnly ". A High Resolution Spatiotemporal Fine Particulate Matter Exposure Assessment Model for the Contiguous United States. "
$para4 = $r4.Paragraphs(1)
$paraStart4 = $para4.Range.Start
$insertRange4 = $d.Range($paraStart4, $paraStart4)
$insertRange4.InsertBefore("Jordan Pennington, Erika Rasnick, Lisa J. Martin, Jocelyn M. Biagini, Tesfaye B. Mersha, Allison Parsons, Gurjit K. Khurana Hershey, Patrick Ryan, ")

[void](FindReplaceOne ". A High Resolution Spatiotemporal Fine Particulate Matter Exposure Assessment Model for the Contiguous United States. " ". Racial Fairness in Precision Medicine: Pediatric Asthma Prediction Algorithms. ")
[void](FindReplaceOne "Environmental Advances" "American Journal of Health Promotion")
[void](FindReplaceOne ". 7:100155. 2022." ". 37(2). 2022.    ")

Write-Host "All edits applied."
